# Update to files: script, anexo.
# Insert a new worksheet "REQ4" (cast of Titanic) before "REQ6", taking the
# first position in the workbook, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Grab a reference to REQ6 to insert the new sheet right before it.
$ws6 = $wb.Worksheets.Item("REQ6")
$new = $wb.Worksheets.Add($ws6)
$new.Name = "REQ4"

# NOTE: after Worksheets.Add() runs, previously-fetched worksheet
# references can go stale (they stop reflecting live formatting), so
# re-fetch every sheet reference we still need by name.
$ws6 = $wb.Worksheets.Item("REQ6")
$new = $wb.Worksheets.Item("REQ4")

# Clone the header cell (style + shared-string slot) from REQ6!A1 ("pelicula")
# onto the new sheet's A1 so the new header reuses the existing gray/centered
# header style instead of minting a new one.
$ws6.Range("A1").Copy($new.Range("A1"))

# Fill in the actor names (Titanic cast) starting at row 2, leaving A1's
# value for last so "actor" lands at the END of the shared-strings table,
# after all the actor names.
$actors = @(
    "Leonardo DiCaprio",
    "Kate Winslet",
    "Billy Zane",
    "Kathy Bates",
    "Frances Fisher",
    "Bernard Hill",
    "Jonathan Hyde",
    "Danny Nucci",
    "David Warner",
    "Bill Paxton",
    "Gloria Stuart",
    "Victor Garber",
    "Suzy Amis"
)

for ($i = 0; $i -lt $actors.Length; $i++) {
    $new.Cells.Item($i + 2, 1).Value = $actors[$i]
}

# Now overwrite the header text (was "pelicula" from the copy above).
$new.Range("A1").Value = "actor"

# Column A width to roughly match the source workbook (20.6328125 chars).
$new.Columns.Item(1).ColumnWidth = 19.83

# Selection shown on the new sheet.
[void]$new.Range("C6").Select()

# Make REQ4 the active tab (it becomes the first, left-most sheet).
$new.Activate()
